# PDH_Hours.xlsx - "added some color highlighting to the PDF hours form"
#
# 1. Re-word the header row so the form is self-explanatory:
#      A1 "Name"  -> "Name to Appear on Certificate"
#      B1 "email" -> "your email address"
#      E1 "Title" -> "Certificate Title"
# 2. Highlight the instructional row (row 2) with a yellow fill so the
#    "fill this in" cells stand out: A2:D2.
# 3. Widen column B so the longer header text is readable.
# 4. Leave the user's selection on E5 (next blank answer cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Header text updates -------------------------------------------------
$ws.Range("A1").Value = "Name to Appear on Certificate"
$ws.Range("E1").Value = "Certificate Title"
$ws.Range("B1").Value = "your email address"

# --- 2. Yellow highlight on row 2 -------------------------------------------
# A2:B2 are blank placeholder cells for the learner to type into.
$ws.Range("A2:B2").Interior.Color = 65535   # RGB(255,255,0)

# C2:D2 already hold the "0 or 1" placeholder text (right aligned) - keep
# that formatting and just add the same yellow fill.
$ws.Range("C2:D2").Interior.Color = 65535   # RGB(255,255,0)

# --- 3. Widen column B for "your email address" -----------------------------
$ws.Range("B1").ColumnWidth = 18.33

# --- 4. Restore the active selection ----------------------------------------
$ws.Range("E5").Select()
